$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: copy cell formatting down from existing template rows ----
# Row 271 carries the "odd" alternating style (14/15/18), row 272 carries the
# "even" alternating style (16/17/19). Row 273 currently carries the special
# "last row" style (20/21) which should move down to the new last row (283).

# preserve the current (last-row) formatting of row 273 for the future row 283
$ws.Range("A273:N273").Copy()
$ws.Range("A283:N283").PasteSpecial(-4122)

# row 273 is no longer the last row -> give it the regular "odd" style
$ws.Range("A271:N271").Copy()
$ws.Range("A273:N273").PasteSpecial(-4122)

# new rows 274-282 alternate even/odd, starting with the "even" style
$ws.Range("A272:N272").Copy()
$ws.Range("A274:N274").PasteSpecial(-4122)
$ws.Range("A271:N271").Copy()
$ws.Range("A275:N275").PasteSpecial(-4122)
$ws.Range("A272:N272").Copy()
$ws.Range("A276:N276").PasteSpecial(-4122)
$ws.Range("A271:N271").Copy()
$ws.Range("A277:N277").PasteSpecial(-4122)
$ws.Range("A272:N272").Copy()
$ws.Range("A278:N278").PasteSpecial(-4122)
$ws.Range("A271:N271").Copy()
$ws.Range("A279:N279").PasteSpecial(-4122)
$ws.Range("A272:N272").Copy()
$ws.Range("A280:N280").PasteSpecial(-4122)
$ws.Range("A271:N271").Copy()
$ws.Range("A281:N281").PasteSpecial(-4122)
$ws.Range("A272:N272").Copy()
$ws.Range("A282:N282").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Step 2: write the actual data (values) row by row, column order B,C,D,E,F,G,H,I,J,K,L,M,N ----
# (this order matches how new shared-string entries were appended in the source diff)

# row 274
$ws.Range("A274").Value2 = 45599.50371961806
$ws.Range("B274").Value = 'leedowon567@naver.com'
$ws.Range("C274").Value = '일본학과'
$ws.Range("D274").Value2 = 20211625.0
$ws.Range("E274").Value = '이도원'
$ws.Range("F274").Value = '대한민국'
$ws.Range("G274").Value = '경제활동인구 / 15세이상 인구'
$ws.Range("H274").Value = '조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자'
$ws.Range("I274").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J274").Value = '평균 : 100만원, 중위값 : 200만원'
$ws.Range("K274").Value = '"19.1%"'
$ws.Range("L274").Value = 'Black'

# row 275
$ws.Range("A275").Value2 = 45599.52109760417
$ws.Range("B275").Value = 'minheart7844@gmail.com'
$ws.Range("C275").Value = '미디어스쿨'
$ws.Range("D275").Value2 = 20242503.0
$ws.Range("E275").Value = '권민'
$ws.Range("F275").Value = '대한민국'
$ws.Range("G275").Value = '취업자 / 경제활동인구'
$ws.Range("H275").Value = '조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자'
$ws.Range("I275").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J275").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K275").Value = '"19.1%"'
$ws.Range("L275").Value = 'Red'
$ws.Range("M275").Value = '나는 사후 장기기증에 참여할 뜻이 없다'

# row 276
$ws.Range("A276").Value2 = 45599.534281782406
$ws.Range("B276").Value = 'ehddn0504@naver.com'
$ws.Range("C276").Value = '소프트웨어학부'
$ws.Range("D276").Value2 = 20245186.0
$ws.Range("E276").Value = '신동우'
$ws.Range("F276").Value = '스페인'
$ws.Range("G276").Value = '취업자 / 15세 이상 인구'
$ws.Range("H276").Value = '조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자'
$ws.Range("I276").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J276").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K276").Value = '"19.1%"'
$ws.Range("L276").Value = 'Black'
$ws.Range("N276").Value = '나는 사후 장기기증에 참여할 뜻이 있다'

# row 277
$ws.Range("A277").Value2 = 45599.54391521991
$ws.Range("B277").Value = 'wlgus4770752@naver.com'
$ws.Range("C277").Value = '의예과'
$ws.Range("D277").Value2 = 20236121.0
$ws.Range("E277").Value = '김지현'
$ws.Range("F277").Value = '대한민국'
$ws.Range("G277").Value = '취업자 / 15세 이상 인구'
$ws.Range("H277").Value = '조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자'
$ws.Range("I277").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("J277").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K277").Value = '"19.1%"'
$ws.Range("L277").Value = 'Red'
$ws.Range("M277").Value = '나는 사후 장기기증에 참여할 뜻이 없다'

# row 278
$ws.Range("A278").Value2 = 45599.54654822916
$ws.Range("B278").Value = 'kjinju0518@naver.com'
$ws.Range("C278").Value = '경영학과'
$ws.Range("D278").Value2 = 20232937.0
$ws.Range("E278").Value = '김진주'
$ws.Range("F278").Value = '대한민국'
$ws.Range("G278").Value = '실업자 / 경제활동인구'
$ws.Range("H278").Value = '조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자'
$ws.Range("I278").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J278").Value = '평균 : 100만원, 중위값 : 200만원'
$ws.Range("K278").Value = '"15%"'
$ws.Range("L278").Value = 'Black'
$ws.Range("N278").Value = '나는 사후 장기기증에 참여할 뜻이 있다'

# row 279
$ws.Range("A279").Value2 = 45599.549399791664
$ws.Range("B279").Value = 'bagj11532@gmail.com'
$ws.Range("C279").Value = '체육학과'
$ws.Range("D279").Value2 = 20244120.0
$ws.Range("E279").Value = '박준형'
$ws.Range("F279").Value = '대한민국'
$ws.Range("G279").Value = '실업자 / 경제활동인구'
$ws.Range("H279").Value = '조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자'
$ws.Range("I279").Value = '평균 : 100만원, 중위값 : 1,000만원'
$ws.Range("J279").Value = '평균 : 100만원, 중위값 : 1,000만원'
$ws.Range("K279").Value = '"25%"'
$ws.Range("L279").Value = 'Red'

# row 280
$ws.Range("A280").Value2 = 45599.551716354166
$ws.Range("B280").Value = 'abcchocoo111@gmail.com'
$ws.Range("C280").Value = '데이터사이언스학부'
$ws.Range("D280").Value2 = 20243206.0
$ws.Range("E280").Value = '김기원'
$ws.Range("F280").Value = '대한민국'
$ws.Range("G280").Value = '취업자 / 15세 이상 인구'
$ws.Range("H280").Value = '조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자'
$ws.Range("I280").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J280").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("K280").Value = '"19.1%"'
$ws.Range("L280").Value = 'Red'
$ws.Range("M280").Value = '나는 사후 장기기증에 참여할 뜻이 없다'

# row 281
$ws.Range("A281").Value2 = 45599.5607603588
$ws.Range("B281").Value = 'ataraxia050508@naver.com'
$ws.Range("C281").Value = '심리학과'
$ws.Range("D281").Value2 = 20242118.0
$ws.Range("E281").Value = '박정호'
$ws.Range("F281").Value = '대한민국'
$ws.Range("G281").Value = '경제활동인구 / 15세이상 인구'
$ws.Range("H281").Value = '조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자'
$ws.Range("I281").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J281").Value = '평균 : 100만원, 중위값 : 1,000만원'
$ws.Range("K281").Value = '"15%"'
$ws.Range("L281").Value = 'Black'
$ws.Range("N281").Value = '나는 사후 장기기증에 참여할 뜻이 있다'

# row 282
$ws.Range("A282").Value2 = 45599.577123449075
$ws.Range("B282").Value = 'sshee718@gmail.com'
$ws.Range("C282").Value = '환경생명공학과'
$ws.Range("D282").Value2 = 20243702.0
$ws.Range("E282").Value = '권도운'
$ws.Range("F282").Value = '대한민국'
$ws.Range("G282").Value = '취업자 / 15세 이상 인구'
$ws.Range("H282").Value = '조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자'
$ws.Range("I282").Value = '평균 : 100만원, 중위값 : 200만원'
$ws.Range("J282").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K282").Value = '"15%"'
$ws.Range("L282").Value = 'Black'

# row 283
$ws.Range("A283").Value2 = 45599.5821505787
$ws.Range("B283").Value = 'dongkyo4@gmail.com'
$ws.Range("C283").Value = '데이터테크'
$ws.Range("D283").Value2 = 20213241.0
$ws.Range("E283").Value = '이동교'
$ws.Range("F283").Value = '대한민국'
$ws.Range("G283").Value = '취업자 / 15세 이상 인구'
$ws.Range("H283").Value = '조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자'
$ws.Range("I283").Value = '평균 : 200만원, 중위값 : 100만원'
$ws.Range("J283").Value = '평균 : 1,000만원, 중위값 : 100만원'
$ws.Range("K283").Value = '"19.1%"'
$ws.Range("L283").Value = 'Red'
